# Insert a new column before column E ("boson") to hold the new "pt_max"
# field, shifting all the existing columns E..O one position to the right
# (F..P), matching the upload's new sheet layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; everything from E onward shifts right.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "pt_max"

# Fill the new column's data rows (2-7) with the value 50.
$ws.Range("E2:E7").Value = 50

# Match the author's final selection state.
[void]$ws.Range("E2:E7").Select()
